# The document contains a standalone paragraph whose entire content is
# the run "This should change". The commit removes that paragraph in
# its entirety (its run AND its own paragraph mark), so the empty
# paragraphs immediately before/after it simply become adjacent.
#
# Deleting Paragraph.Range (as opposed to deleting just the found text
# range) removes the paragraph's mark along with its text, which is
# what collapses the paragraph out of the document instead of leaving
# an empty paragraph behind.

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*This should change*") {
        $para.Range.Delete()
        break
    }
}
